$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19455
$ws.Range("B3").Value = 14562
$ws.Range("B4").Value = 1816
$ws.Range("B5").Value = 18019
